$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (row index => column letter => value)
# Columns: E, G, H, I, J, K, M, N, O, P, Q, R, S, T

$rowsData = @{
    2 = @{
        E = 3
        G = 1.298902666666667
        H = 3.896708
        I = 0.7196693520699016
        J = 0.7196693520699017
        K = 3
        M = 8.922308333333334
        N = 26.766925
        O = 0.6138261687668722
        P = 0.6138261687668722
        Q = 11.58921008698889
        R = 104.3028907829
        S = 0.441751881160005
        T = 0.4417518811600051
    }
    3 = @{
        E = 3
        G = 1.298902666666667
        H = 3.896708
        I = 0.7196693520699016
        J = 0.7196693520699017
        K = 3
        M = 5.613253666666666
        N = 16.839761
        O = 0.3861738312331279
        P = 0.3861738312331279
        Q = 7.291070156309778
        R = 65.619631406788
        S = 0.2779174709098967
        T = 0.2779174709098967
    }
    4 = @{
        E = 3
        G = 0.5059576666666666
        H = 1.517873
        I = 0.2803306479300983
        J = 0.2803306479300984
        K = 3
        M = 8.922308333333334
        N = 26.766925
        O = 0.6138261687668722
        P = 0.6138261687668722
        Q = 4.514310305613889
        R = 40.628792750525
        S = 0.1720742876068672
        T = 0.1720742876068672
    }
    5 = @{
        E = 3
        G = 0.5059576666666666
        H = 1.517873
        I = 0.2803306479300983
        J = 0.2803306479300984
        K = 3
        M = 5.613253666666666
        N = 16.839761
        O = 0.3861738312331279
        P = 0.3861738312331279
        Q = 2.840068727594777
        R = 25.560618548353
        S = 0.1082563603232312
        T = 0.1082563603232312
    }
}

foreach ($rowNum in $rowsData.Keys) {
    $cols = $rowsData[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
